$d = $word.ActiveDocument

# --- Locate the target paragraph ("this is a piece of code to do simple additon") ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "this is a piece of code to do simple additon*") {
        $target = $cand
        break
    }
}

# --- Split the misspelled word "additon" into its own run while fixing the typo ---
# "this is a piece of code to do simple " is 37 characters (incl. trailing space)
$prefix = "this is a piece of code to do simple "
$paraStart = $target.Range.Start
$splitPos = $paraStart + $prefix.Length
$paraEnd = $target.Range.End

# Range covering the old, misspelled tail ("additon"), excluding the paragraph mark
$tail = $d.Range($splitPos, $paraEnd - 1)

# Nudge formatting so this substring is emitted as its own run, then restore it,
# and fix the typo while we're at it.
$tail.Bold = 1
$tail.Text = "addition"
$tail.Bold = 0

# --- Add a new paragraph after it with the extra line of text ---
$target = $d.Paragraphs.Item($target.Index)
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($target.Index + 1)
$newPara.Range.Text = "this code is written by python language"

Write-Output "done"
